$d = $word.ActiveDocument
$r = $d.Content
$r.Start = 0
$r.End = 0
Write-Output $r.WordOpenXML.Substring(0, 500)
